# Generate Report for Handback
# Update the "Latest Handback DateTime" and "Error Detail" values for the
# 1d810fc8-2299-46f8-a074-2c75f9035a63 row on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 4 corresponds to 1d810fc8-2299-46f8-a074-2c75f9035a63 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("L4").Value = "2017-02-21 11:07:13"
$wsZhCn.Range("R4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/f0fa68b2cedb0678dec8cfee90fe35ec3280b44c/e2e/1d810fc8-2299-46f8-a074-2c75f9035a63.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/2f2cc2edf61b9592c5a5679de992c9f8255c2241/e2e/1d810fc8-2299-46f8-a074-2c75f9035a63.md."

# --- de-de sheet: row 4 corresponds to 1d810fc8-2299-46f8-a074-2c75f9035a63 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("L4").Value = "2017-02-21 11:07:35"
